$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new species/model name text, new relative-abundance value)
# Text values already include the literal leading/trailing apostrophe that is
# part of the target data (e.g. "'Bacteroides_..._DSM_14838.mat'"). A leading
# apostrophe typed into Excel is normally consumed as a "force text" marker,
# so it is doubled here ("''...") to make Excel keep a single literal
# apostrophe at the start of the stored string.
$rows = @(
  @{ Row = 2;  Name = "'Bacteroides_cellulosilyticus_DSM_14838.mat'"; Value = 0.016 },
  @{ Row = 3;  Name = "'Bacteroides_coprocola_M16_DSM_17136.mat'";    Value = 0.003 },
  @{ Row = 4;  Name = "'Bacteroides_coprophilus_DSM_18228.mat'";      Value = 0.047 },
  @{ Row = 5;  Name = "'Bacteroides_fluxus_YIT_12057.mat'";           Value = 0.065 },
  @{ Row = 6;  Name = "'Bacteroides_oleiciplenus_YIT_12058.mat'";     Value = 0.013 },
  @{ Row = 7;  Name = "'Bacteroides_ovatus_ATCC_8483.mat'";           Value = 0.231 },
  @{ Row = 8;  Name = "'Bacteroides_salyersiae_WAL_10018.mat'";       Value = 0.305 },
  @{ Row = 9;  Name = "'Bacteroides_stercoris_ATCC_43183.mat'";       Value = -0 },
  @{ Row = 10; Name = "'Bacteroides_thetaiotaomicron_VPI_5482.mat'";  Value = -0 },
  @{ Row = 11; Name = "'Bacteroides_uniformis_ATCC_8492.mat'";        Value = -0 },
  @{ Row = 12; Name = "'Bacteroides_vulgatus_ATCC_8482.mat'";         Value = 0.321 },
  @{ Row = 13; Name = "'Bifidobacterium_animalis_lactis_AD011.mat'";  Value = 0 },
  @{ Row = 14; Name = "'Enterococcus_faecalis_OG1RF_ATCC_47077.mat'"; Value = -0 },
  @{ Row = 15; Name = "'Flavonifractor_plautii_ATCC_29863.mat'";      Value = -0 },
  @{ Row = 16; Name = "'Lactobacillus_plantarum_JDM1.mat'";           Value = 0 },
  @{ Row = 17; Name = "'Odoribacter_laneus_YIT_12061.mat'";           Value = -0 },
  @{ Row = 18; Name = "'Parabacteroides_distasonis_ATCC_8503.mat'";   Value = 0 },
  @{ Row = 19; Name = "'Parabacteroides_johnsonii_DSM_18315.mat'";    Value = 0 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = "'" + $item.Name
    # Clear the automatic "quote prefix" formatting Excel applies when a
    # cell's text begins with an apostrophe, so only the cell content (not
    # its style) changes.
    $ws.Range("B$r").Style = "Normal"
    $ws.Range("C$r").Value = $item.Value
}
